# Remove the "preencoded.png" watermark picture (linked to https://gamma.app)
# that appears near the bottom-right corner of every content slide except
# the title slide. Identify it by its fixed position/size (EMU):
#   off  x=12242153 y=7589520
#   ext  cx=2296807 cy=548640
# which in points is Left=963.949 Top=597.6 Width=180.851 Height=43.2.

$p = $ppt.ActivePresentation

$targetLeft = 963.949
$targetTop = 597.6
$targetWidth = 180.851
$targetHeight = 43.2
$tolerance = 0.5

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $s = $p.Slides.Item($slideIdx)
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Type -eq 13) {
            $dLeft = [Math]::Abs($sh.Left - $targetLeft)
            $dTop = [Math]::Abs($sh.Top - $targetTop)
            $dWidth = [Math]::Abs($sh.Width - $targetWidth)
            $dHeight = [Math]::Abs($sh.Height - $targetHeight)
            if (($dLeft -lt $tolerance) -and ($dTop -lt $tolerance) -and ($dWidth -lt $tolerance) -and ($dHeight -lt $tolerance)) {
                $sh.Delete()
            }
        }
    }
}
